# Fruta / hortaliza, semanal
# Insert a new weekly price block (3 rows) for "Terminal La Palmera de La Serena - Piña"
# above the existing block that starts at row 546, shifting all following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 546:548 (existing rows 546:570 shift down to 549:573)
$ws.Rows("546:548").Insert()

# Row 546 - Especial / $/caja 10 unidades
$ws.Cells.Item(546, 1).Value = 8
$ws.Cells.Item(546, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(546, 3).Value = "Coquimbo"
$ws.Cells.Item(546, 4).Value = 44509
$ws.Cells.Item(546, 5).Value = 4
$ws.Cells.Item(546, 6).Value = "Fruta"
$ws.Cells.Item(546, 7).Value = 100108
$ws.Cells.Item(546, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(546, 9).Value = 100108005
$ws.Cells.Item(546, 10).Value = "Piña"
$ws.Cells.Item(546, 11).Value = "Caramelo"
$ws.Cells.Item(546, 12).Value = "Especial"
$ws.Cells.Item(546, 13).Value = 432
$ws.Cells.Item(546, 14).Value = 19000
$ws.Cells.Item(546, 15).Value = 20000
$ws.Cells.Item(546, 16).Value = 19500
$ws.Cells.Item(546, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(546, 18).Value = "Ecuador"
$ws.Cells.Item(546, 19).Value = 1950
$ws.Cells.Item(546, 20).Value = 10

# Row 547 - Primera / $/caja 12 unidades
$ws.Cells.Item(547, 1).Value = 8
$ws.Cells.Item(547, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(547, 3).Value = "Coquimbo"
$ws.Cells.Item(547, 4).Value = 44509
$ws.Cells.Item(547, 5).Value = 4
$ws.Cells.Item(547, 6).Value = "Fruta"
$ws.Cells.Item(547, 7).Value = 100108
$ws.Cells.Item(547, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(547, 9).Value = 100108005
$ws.Cells.Item(547, 10).Value = "Piña"
$ws.Cells.Item(547, 11).Value = "Caramelo"
$ws.Cells.Item(547, 12).Value = "Primera"
$ws.Cells.Item(547, 13).Value = 432
$ws.Cells.Item(547, 14).Value = 19000
$ws.Cells.Item(547, 15).Value = 20000
$ws.Cells.Item(547, 16).Value = 19500
$ws.Cells.Item(547, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(547, 18).Value = "Ecuador"
$ws.Cells.Item(547, 19).Value = 1625
$ws.Cells.Item(547, 20).Value = 12

# Row 548 - Segunda / $/caja 14 unidades
$ws.Cells.Item(548, 1).Value = 8
$ws.Cells.Item(548, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(548, 3).Value = "Coquimbo"
$ws.Cells.Item(548, 4).Value = 44509
$ws.Cells.Item(548, 5).Value = 4
$ws.Cells.Item(548, 6).Value = "Fruta"
$ws.Cells.Item(548, 7).Value = 100108
$ws.Cells.Item(548, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(548, 9).Value = 100108005
$ws.Cells.Item(548, 10).Value = "Piña"
$ws.Cells.Item(548, 11).Value = "Caramelo"
$ws.Cells.Item(548, 12).Value = "Segunda"
$ws.Cells.Item(548, 13).Value = 432
$ws.Cells.Item(548, 14).Value = 19000
$ws.Cells.Item(548, 15).Value = 20000
$ws.Cells.Item(548, 16).Value = 19500
$ws.Cells.Item(548, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(548, 18).Value = "Ecuador"
$ws.Cells.Item(548, 19).Value = 1393
$ws.Cells.Item(548, 20).Value = 14
